# Insert a new row at position 104 (shifts existing rows 104-145 down to 105-146)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(104).Insert()

# Populate the newly inserted row 104 with the new data record
$ws.Range("A104").Value = 7
$ws.Range("B104").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C104").Value = "Ñuble"
$ws.Range("D104").Value = 44988
$ws.Range("E104").Value = 16
$ws.Range("F104").Value = 100112021
$ws.Range("G104").Value = "Ají"
$ws.Range("H104").Value = "Cacho cabra verde"
$ws.Range("I104").Value = "Primera"
$ws.Range("J104").Value = 30
$ws.Range("K104").Value = 14000
$ws.Range("L104").Value = 14000
$ws.Range("M104").Value = 14000
$ws.Range("N104").Value = "$/saco 25 kilos"
$ws.Range("O104").Value = "Región del Maule"
$ws.Range("P104").Value = 560
$ws.Range("Q104").Value = 25
$ws.Range("R104").Value = "Hortaliza"
